## Add a new worksheet "2022_7" at the end of the workbook with user-log data
## (mirrors the structure of the existing 2022_x sheets).

$wb = $excel.ActiveWorkbook

# Insert the new worksheet after the last existing sheet so it lands at the end.
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "2022_7"

# Helper-free direct cell writes.  For values that look numeric (phone
# numbers, numeric-strings such as "1") we force a text number format
# first so Excel keeps them as text instead of silently coercing them to
# numbers.

function Set-TextCell($range, [string]$value) {
    $range.NumberFormat = "@"
    $range.Value = $value
}

# ---- Header row ----
$ws.Range("A1").Value = "date"
$ws.Range("B1").Value = "time"
$ws.Range("C1").Value = "phoneNumber"
$ws.Range("D1").Value = "model"
$ws.Range("E1").Value = "path"
$ws.Range("F1").Value = "action"
$ws.Range("G1").Value = "status"
$ws.Range("H1").Value = "description"
$ws.Range("I1").Value = "failureReason"
$ws.Range("J1").Value = "userId"
$ws.Range("K1").Value = "modelId"
$ws.Range("L1").Value = "lastName"
$ws.Range("M1").Value = "firstName"

# ---- Row 2 ----
$ws.Range("A2").Value = "Sat Jul 09 2022"
$ws.Range("B2").Value = "16:12:36 GMT+0000 (Coordinated Universal Time)"
Set-TextCell $ws.Range("C2") "+22892942601"
$ws.Range("D2").Value = "User"
$ws.Range("E2").Value = "/api/auth/send-otp"
$ws.Range("F2").Value = "request"
$ws.Range("G2").Value = "failed"
$ws.Range("H2").Value = "+22892942601 request to receive otp"
$ws.Range("I2").Value = "error.invalid"
$ws.Range("J2").Value = "+22892942601 request to receive otp"
$ws.Range("K2").Value = "error.invalid"

# ---- Row 3 ----
$ws.Range("A3").Value = "Sat Jul 09 2022"
$ws.Range("B3").Value = "16:13:00 GMT+0000 (Coordinated Universal Time)"
Set-TextCell $ws.Range("C3") "22892942601"
$ws.Range("D3").Value = "User"
$ws.Range("E3").Value = "/api/auth/send-otp"
$ws.Range("F3").Value = "request"
$ws.Range("G3").Value = "failed"
$ws.Range("H3").Value = "22892942601 request to receive otp"
$ws.Range("I3").Value = "error.userNotFound"

# ---- Row 4 ----
$ws.Range("A4").Value = "Sat Jul 09 2022"
$ws.Range("B4").Value = "16:13:52 GMT+0000 (Coordinated Universal Time)"
Set-TextCell $ws.Range("C4") "22892942601"
$ws.Range("D4").Value = "User"
$ws.Range("E4").Value = "/api/auth/send-otp"
$ws.Range("F4").Value = "request"
$ws.Range("G4").Value = "failed"
$ws.Range("H4").Value = "22892942601 request to receive otp"
$ws.Range("I4").Value = "error.userNotFound"

# ---- Row 5 ----
$ws.Range("A5").Value = "Sat Jul 09 2022"
$ws.Range("B5").Value = "16:17:50 GMT+0000 (Coordinated Universal Time)"
Set-TextCell $ws.Range("C5") "22892942601"
$ws.Range("D5").Value = "User"
$ws.Range("E5").Value = "/api/auth/send-otp"
$ws.Range("F5").Value = "request"
$ws.Range("G5").Value = "failed"
$ws.Range("H5").Value = "22892942601 request to receive otp"
$ws.Range("I5").Value = "getaddrinfo EAI_AGAIN dashboard.smszedekaa.com"

# ---- Row 6 ----
$ws.Range("A6").Value = "Sat Jul 09 2022"
$ws.Range("B6").Value = "16:27:47 GMT+0000 (Coordinated Universal Time)"
Set-TextCell $ws.Range("C6") "22892942601"
$ws.Range("D6").Value = "User"
$ws.Range("E6").Value = "/api/auth/send-otp"
$ws.Range("F6").Value = "request"
$ws.Range("G6").Value = "succeeded"
$ws.Range("H6").Value = "22892942601 request to receive otp"

# ---- Row 7 ----
$ws.Range("A7").Value = "Sat Jul 09 2022"
$ws.Range("B7").Value = "16:28:54 GMT+0000 (Coordinated Universal Time)"
Set-TextCell $ws.Range("C7") "22892942601"
$ws.Range("D7").Value = "User"
$ws.Range("E7").Value = "/api/auth/verify-otp"
$ws.Range("F7").Value = "request"
$ws.Range("G7").Value = "succeeded"
$ws.Range("H7").Value = "22892942601 request to receive otp"
$ws.Range("J7").Value = 1
Set-TextCell $ws.Range("K7") "1"

# ---- Row 8 ----
$ws.Range("A8").Value = "Sat Jul 09 2022"
$ws.Range("B8").Value = "16:29:55 GMT+0000 (Coordinated Universal Time)"
$ws.Range("D8").Value = "User"
$ws.Range("E8").Value = "/api/auth/complete-infos"
$ws.Range("F8").Value = "edit"
$ws.Range("G8").Value = "succeeded"
$ws.Range("H8").Value = "    edit his infos"
$ws.Range("J8").Value = 1
Set-TextCell $ws.Range("K8") "1"
Set-TextCell $ws.Range("L8") ""
Set-TextCell $ws.Range("M8") ""

# ---- Row 9 ----
$ws.Range("A9").Value = "Sun Jul 10 2022"
$ws.Range("B9").Value = "01:27:53 GMT+0000 (Coordinated Universal Time)"
Set-TextCell $ws.Range("C9") "22892942601"
$ws.Range("D9").Value = "User"
$ws.Range("E9").Value = "/api/auth/send-otp"
$ws.Range("F9").Value = "request"
$ws.Range("G9").Value = "succeeded"
$ws.Range("H9").Value = "22892942601 request to receive otp"

# ---- Row 10 ----
$ws.Range("A10").Value = "Sun Jul 10 2022"
$ws.Range("B10").Value = "01:28:20 GMT+0000 (Coordinated Universal Time)"
Set-TextCell $ws.Range("C10") "22892942601"
$ws.Range("D10").Value = "User"
$ws.Range("E10").Value = "/api/auth/verify-otp"
$ws.Range("F10").Value = "request"
$ws.Range("G10").Value = "succeeded"
$ws.Range("H10").Value = "22892942601 request to receive otp"
$ws.Range("J10").Value = 1
Set-TextCell $ws.Range("K10") "1"

# Avoid leaving the new sheet "active" (the original sheets carry no
# tabSelected/selection markup), mirroring the sibling sheets' plain view.
$wb.Worksheets.Item(1).Activate()
